# Fix the CasesTab ("startup"!B2) Cypher query: the `Cohort` coalesce line
# referenced a variable (`co`) that is no longer wired into the query result
# set the way the rest of the workbook expects, so drop the trailing
# ", coalesce(co.cohort_description, '') AS `Cohort`" line (and the now
# dangling trailing blank line) from the stored query text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesTabQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['West Highland White Terrier'] `nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesTabQuery

# Restore the selection to the cell that was actually edited.
$ws.Range("B2").Select()
